# Append a new data row (row 65) to the Adafruit IO export sheet,
# mirroring the existing rows of timestamp/feed/value/lat/long/elevation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 65

# Column C holds a numeric-looking value ("25") that must stay text (as all
# the other data in this sheet is stored as text), so force a text number
# format on that cell before writing the value.
$ws.Range("C$newRow").NumberFormat = "@"

$ws.Range("A$newRow").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$newRow").Value = "temperature"
$ws.Range("C$newRow").Value = "25"
$ws.Range("D$newRow").Value = "N/A"
$ws.Range("E$newRow").Value = "N/A"
$ws.Range("F$newRow").Value = "N/A"
